$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update Date value in B8 ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-05-02T17:24:18+00:00"

# --- Sheet "Elements": update Type(s) and Short for AuthorPerson.1 row (row 3) ---
$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("K3").Value = "Reference(https://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/PSIdNat|https://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/MatriculeINS|https://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/SNR|https://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/IdentifiantSysteme)`n"
$wsElem.Range("L3").Value = "Identifiant"
$wsElem.Range("M3").Value = "Identifiant"

# --- Column K width on Elements sheet ---
# NOTE: Excel's ColumnWidth (character units) is converted to the stored
# OOXML "width" (pixel/MDW grid) with rounding, so 254.15 is the
# ColumnWidth input that rounds to the target stored width of 255.0.
$wsElem.Columns.Item(11).ColumnWidth = 254.15
